$wb = $excel.ActiveWorkbook

# ===================================================================
# This workbook reports localization handoff status for two files:
#   4cccf715-1532-4194-b4da-fc9bb09035d5.md
#   5b587bc9-c77c-4c85-aed3-aaa9de3d452a.md
# A new handoff report is generated: 5b587bc9... is now listed first
# (still "In Translation"), and 4cccf715... has just been generated /
# handed off, moving to "Ready for handoff" with a fresh handoff time.
# ===================================================================

# ----- Sheet: Overview -----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A1:D3").ClearContents()
$wsOverview.Hyperlinks.Delete()

$wsOverview.Range("A1").Value = "File Name"
$wsOverview.Range("B1").Value = "zh-cn"
$wsOverview.Range("C1").Value = "de-de"
$wsOverview.Range("D1").Value = "Latest Handoff Date"

$wsOverview.Range("A2").Value = "5b587bc9-c77c-4c85-aed3-aaa9de3d452a.md"
$wsOverview.Range("B2").Value = "In Translation"
$wsOverview.Range("C2").Value = "In Translation"
$wsOverview.Range("D2").Value = "2016-03-18 04:03:04"

$wsOverview.Range("A3").Value = "4cccf715-1532-4194-b4da-fc9bb09035d5.md"
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-03-18 04:03:42"

$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/2ce1a1ee5097d1d34f75206e203d154bef4bf93e/e2e/4cccf715-1532-4194-b4da-fc9bb09035d5.md", "", "", "5b587bc9-c77c-4c85-aed3-aaa9de3d452a.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/2ce1a1ee5097d1d34f75206e203d154bef4bf93e/e2e/5b587bc9-c77c-4c85-aed3-aaa9de3d452a.md", "", "", "4cccf715-1532-4194-b4da-fc9bb09035d5.md")

# ----- Sheet: zh-cn -----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A1:K3").ClearContents()
$wsZh.Hyperlinks.Delete()

$wsZh.Range("A1").Value = "Source File Name"
$wsZh.Range("B1").Value = "File Extension"
$wsZh.Range("C1").Value = "Status"
$wsZh.Range("D1").Value = "Latest Handoff File"
$wsZh.Range("E1").Value = "Latest Handoff Datetime"
$wsZh.Range("F1").Value = "Latest Target File"
$wsZh.Range("G1").Value = "Latest Handback File"
$wsZh.Range("H1").Value = "Latest Handback DateTime"
$wsZh.Range("I1").Value = "Handoff Reason"
$wsZh.Range("J1").Value = "Dependency From"
$wsZh.Range("K1").Value = "Error Detail"

$wsZh.Range("A2").Value = "5b587bc9-c77c-4c85-aed3-aaa9de3d452a.md"
$wsZh.Range("B2").Value = ".md"
$wsZh.Range("C2").Value = "In Translation"
$wsZh.Range("D2").Value = "5b587bc9-c77c-4c85-aed3-aaa9de3d452a.632735bdbcaaf6d8835780987cce7e370fe186ed.zh-cn.xlf"
$wsZh.Range("E2").Value = "2016-03-18 04:03:00"
$wsZh.Range("H2").Value = "0001-01-01 00:00:00"
$wsZh.Range("I2").Value = "Include"

$wsZh.Range("A3").Value = "4cccf715-1532-4194-b4da-fc9bb09035d5.md"
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "4cccf715-1532-4194-b4da-fc9bb09035d5.915e06189cd0a3f3b3a05cc0384616ccc1c16fa3.zh-cn.xlf"
$wsZh.Range("E3").Value = "2016-03-18 04:03:40"
$wsZh.Range("H3").Value = "0001-01-01 00:00:00"
$wsZh.Range("I3").Value = "Include"

$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/2ce1a1ee5097d1d34f75206e203d154bef4bf93e/e2e/4cccf715-1532-4194-b4da-fc9bb09035d5.md", "", "", "5b587bc9-c77c-4c85-aed3-aaa9de3d452a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/2ce1a1ee5097d1d34f75206e203d154bef4bf93e/e2e/4cccf715-1532-4194-b4da-fc9bb09035d5.md", "", "", ".md")
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/74290a5681f5d610557332750c241e4e33cac69a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4cccf715-1532-4194-b4da-fc9bb09035d5.915e06189cd0a3f3b3a05cc0384616ccc1c16fa3.zh-cn.xlf", "", "", "5b587bc9-c77c-4c85-aed3-aaa9de3d452a.632735bdbcaaf6d8835780987cce7e370fe186ed.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/2ce1a1ee5097d1d34f75206e203d154bef4bf93e/e2e/5b587bc9-c77c-4c85-aed3-aaa9de3d452a.md", "", "", "4cccf715-1532-4194-b4da-fc9bb09035d5.md")
$wsZh.Hyperlinks.Add($wsZh.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/2ce1a1ee5097d1d34f75206e203d154bef4bf93e/e2e/5b587bc9-c77c-4c85-aed3-aaa9de3d452a.md", "", "", ".md")
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/74290a5681f5d610557332750c241e4e33cac69a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/5b587bc9-c77c-4c85-aed3-aaa9de3d452a.632735bdbcaaf6d8835780987cce7e370fe186ed.zh-cn.xlf", "", "", "4cccf715-1532-4194-b4da-fc9bb09035d5.915e06189cd0a3f3b3a05cc0384616ccc1c16fa3.zh-cn.xlf")

# ----- Sheet: de-de -----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A1:K3").ClearContents()
$wsDe.Hyperlinks.Delete()

$wsDe.Range("A1").Value = "Source File Name"
$wsDe.Range("B1").Value = "File Extension"
$wsDe.Range("C1").Value = "Status"
$wsDe.Range("D1").Value = "Latest Handoff File"
$wsDe.Range("E1").Value = "Latest Handoff Datetime"
$wsDe.Range("F1").Value = "Latest Target File"
$wsDe.Range("G1").Value = "Latest Handback File"
$wsDe.Range("H1").Value = "Latest Handback DateTime"
$wsDe.Range("I1").Value = "Handoff Reason"
$wsDe.Range("J1").Value = "Dependency From"
$wsDe.Range("K1").Value = "Error Detail"

$wsDe.Range("A2").Value = "5b587bc9-c77c-4c85-aed3-aaa9de3d452a.md"
$wsDe.Range("B2").Value = ".md"
$wsDe.Range("C2").Value = "In Translation"
$wsDe.Range("D2").Value = "5b587bc9-c77c-4c85-aed3-aaa9de3d452a.632735bdbcaaf6d8835780987cce7e370fe186ed.de-de.xlf"
$wsDe.Range("E2").Value = "2016-03-18 04:03:04"
$wsDe.Range("H2").Value = "0001-01-01 00:00:00"
$wsDe.Range("I2").Value = "Include"

$wsDe.Range("A3").Value = "4cccf715-1532-4194-b4da-fc9bb09035d5.md"
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "4cccf715-1532-4194-b4da-fc9bb09035d5.915e06189cd0a3f3b3a05cc0384616ccc1c16fa3.de-de.xlf"
$wsDe.Range("E3").Value = "2016-03-18 04:03:42"
$wsDe.Range("H3").Value = "0001-01-01 00:00:00"
$wsDe.Range("I3").Value = "Include"

$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/2ce1a1ee5097d1d34f75206e203d154bef4bf93e/e2e/4cccf715-1532-4194-b4da-fc9bb09035d5.md", "", "", "5b587bc9-c77c-4c85-aed3-aaa9de3d452a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/2ce1a1ee5097d1d34f75206e203d154bef4bf93e/e2e/4cccf715-1532-4194-b4da-fc9bb09035d5.md", "", "", ".md")
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cba815969b818987d14a12747b67e202d9c5fc5a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4cccf715-1532-4194-b4da-fc9bb09035d5.915e06189cd0a3f3b3a05cc0384616ccc1c16fa3.de-de.xlf", "", "", "5b587bc9-c77c-4c85-aed3-aaa9de3d452a.632735bdbcaaf6d8835780987cce7e370fe186ed.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/2ce1a1ee5097d1d34f75206e203d154bef4bf93e/e2e/5b587bc9-c77c-4c85-aed3-aaa9de3d452a.md", "", "", "4cccf715-1532-4194-b4da-fc9bb09035d5.md")
$wsDe.Hyperlinks.Add($wsDe.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/2ce1a1ee5097d1d34f75206e203d154bef4bf93e/e2e/5b587bc9-c77c-4c85-aed3-aaa9de3d452a.md", "", "", ".md")
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cba815969b818987d14a12747b67e202d9c5fc5a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/5b587bc9-c77c-4c85-aed3-aaa9de3d452a.632735bdbcaaf6d8835780987cce7e370fe186ed.de-de.xlf", "", "", "4cccf715-1532-4194-b4da-fc9bb09035d5.915e06189cd0a3f3b3a05cc0384616ccc1c16fa3.de-de.xlf")
